$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data update (GitHub Actions symbol list refresh)
# Each Price (D) / Volume(1h) (E) cell is stored as text in the source data,
# so we force text formatting before writing to avoid Excel auto-converting
# these into numeric/percentage values.
$cells = @{
    'D2' = '278.77'
    'E2' = '6.59%'
    'E3' = '2.92%'
    'D4' = '4.810'
    'D5' = '0.06332'
    'E5' = '2.44%'
    'D6' = '6.949'
    'E6' = '3.65%'
    'D7' = '3.408'
    'E7' = '7.34%'
    'D8' = '0.8801'
    'E8' = '3.51%'
    'D9' = '0.9581'
    'E9' = '5.13%'
    'D10' = '0.1472'
    'E10' = '4.91%'
    'D11' = '0.05133'
    'E11' = '-0.29%'
    'D12' = '0.07324'
    'E12' = '3.15%'
    'D13' = '0.03133'
    'E13' = '0.86%'
    'D14' = '0.09071'
    'E14' = '0.28%'
    'D15' = '0.001569'
    'E15' = '2.53%'
    'D16' = '0.0006263'
    'E16' = '1.25%'
    'D17' = '0.005773'
    'E17' = '-3.59%'
    'D18' = '3.459'
    'E18' = '0.30%'
    'D19' = '2.295'
    'E19' = '4.93%'
    'D20' = '0.3140'
    'E20' = '2.19%'
    'E21' = '0.06%'
    'D22' = '3.863'
    'E22' = '-5.38%'
    'D23' = '0.04327'
    'E23' = '2.22%'
    'D24' = '0.001179'
    'E24' = '0.27%'
    'D25' = '0.004289'
    'E25' = '5.94%'
    'D26' = '0.0001201'
    'E26' = '-0.02%'
    'D27' = '0.0001691'
    'E27' = '3.12%'
    'D40' = '0.04089'
    'E40' = '2.45%'
    'D41' = '0.006769'
    'E41' = '63.47%'
    'D42' = '0.1162'
    'E42' = '4.40%'
    'D43' = '0.002201'
    'E43' = '2.73%'
    'D44' = '0.01291'
    'E44' = '-2.82%'
    'D45' = '0.00005210'
    'E45' = '0.90%'
    'E46' = '-0.02%'
    'E47' = '843.34%'
    'D48' = '0.02251'
    'E48' = '-33.83%'
    'D49' = '0.00002101'
    'E49' = '-0.02%'
    'D50' = '0.0002001'
    'E50' = '-0.02%'
}

foreach ($addr in $cells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cells[$addr]
    $cell.Style = "Normal"
}
